$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.156.03"
$ws.Range("E2").Value = "  -2.36%  "

$ws.Range("D3").Value = "1.839.38"
$ws.Range("E3").Value = "  -1.55%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.9999"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.03%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "239.97"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -2.81%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.6817"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -2.75%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.000"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.04%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2999"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -2.93%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07461"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -4.18%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "23.24"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -2.92%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07636"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -2.74%  "

$ws.Range("D12").Value = "1.842.89"
$ws.Range("E12").Value = "  -1.24%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "5.045"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -2.86%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.6810"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -2.30%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "87.95"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -5.33%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "6.102"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -8.13%  "

$ws.Range("D17").Value = "29.175.73"
$ws.Range("E17").Value = "  -2.24%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.000008198"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -2.67%  "

$ws.Range("D19").Value = "2.091.48"
$ws.Range("E19").Value = "  -0.83%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "230.07"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -5.89%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "12.53"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -2.57%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.9994"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.05%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.359"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -4.20%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "1.001"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.06%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "160.25"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.00%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.1441"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -4.88%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "8.703"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -3.22%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "18.10"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -1.94%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.502"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -2.71%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "4.265"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -0.73%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.141"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -2.46%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.193"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -1.00%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.05326"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +4.13%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.7541"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -4.55%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.855"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -4.44%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.133"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -3.03%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.688"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.91%  "

$ws.Range("D38").Value = "1.312.22"
$ws.Range("E38").Value = "  -1.45%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01829"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -3.16%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.723"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -0.91%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.9464"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -0.84%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "6.011"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -1.07%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "104.57"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -2.70%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.9992"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.06%  "

$ws.Range("D45").Value = "1.992.20"
$ws.Range("E45").Value = "  -1.01%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.5187"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.30%  "

$ws.Range("B47").Value = "BabyDogeCoin"
$ws.Range("C47").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.00000000122"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -3.10%  "

$ws.Range("B48").Value = "Aave"
$ws.Range("C48").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "64.48"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -1.66%  "

$ws.Range("B49").Value = "EnergySwap"
$ws.Range("C49").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "9.489"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -3.68%  "

$ws.Range("B50").Value = "RenderToken"
$ws.Range("C50").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.772"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -1.68%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.07641"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +16.19%  "
